$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    0.084042413456785425,
    -0.0059999999395081716,
    -0.0039999999532653874,
    -0.007999999912545519,
    -0.0029999999615029083,
    -0.0019999999711561856,
    -0.0099999998926678657,
    -0.0099999998931825651,
    -0.0019999999755455633,
    -0.001999999978844258,
    -0.0029999999691732171,
    -0.0034999999652836067,
    -0.0034999999721669894,
    0.03445413421544874,
    -0.00099999999928890304,
    -0.0019999999887021502,
    -0.0019999999878930197,
    -0.003999999967744472,
    -0.0039999999593334223,
    0.0086011840801081263,
    -0.0039999999601318947,
    0.0020624633656112223,
    -0.004999999939853339,
    -0.019999999784526601,
    -0.019999999781303401,
    -0.0024999999661829975,
    -0.002499999966200761,
    -0.0019999999717921213,
    -0.0069999999228764764,
    -0.059999999397998049,
    0.069486373468924612,
    -0.0099999999072615253,
    0.010122491249346055
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $values[$i]
}
